# Workbook has two sheets:
#   1) addMultipleCustomerTest  (currently the active/selected tab)
#   2) openAccountMethod
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update the two mislabeled cells on the first sheet ---
# A4 was "Keyword" -> "No"
# B4 was "Driven"  -> "Cust"
$ws1.Range("A4").Value = "No"
$ws1.Range("B4").Value = "Cust"

# --- Move the selection on sheet 1 down to A5 ---
$ws1.Range("A5").Select()

# --- Switch the active tab to the second sheet (openAccountMethod) ---
$ws2.Activate()

Write-Output "done"
